# First pass of new method
#
# For each "~UC_Sets: T_xx:" header block in the constraint sheets, the
# label that used to live in column B has moved to column A, and on the
# row below it the two values that used to sit in A/B (the "~UC_T" /
# "~TFM_INS" / "~TFM_UPD" marker and the "~UC_Sets: R_x: AllRegions"
# marker) have swapped columns as well:
#
#   Row N   :  B{N}   "~UC_Sets: T_xx: "      -> A{N}   "~UC_Sets: T_xx: "
#   Row N+1 :  A{N+1} "~UC_T"/"~TFM_INS"/...  -> B{N+1} "~UC_T"/"~TFM_INS"/...
#              B{N+1} "~UC_Sets: R_x: ..."    -> A{N+1} "~UC_Sets: R_x: ..."
#
# This happens on every sheet of the workbook, at every occurrence of the
# pattern. Apply it generically by scanning each sheet's used range for
# cells in column B whose text starts with "~UC_Sets: T_" (the row-N
# marker), then swapping the two cells in the row directly below it.

$wb = $excel.ActiveWorkbook

function Swap-UcSetsBlock($ws, $rowTitle) {
    $rowHeader = $rowTitle + 1

    $titleCellB = $ws.Cells.Item($rowTitle, 2)
    $titleText = $titleCellB.Formula

    $headerCellA = $ws.Cells.Item($rowHeader, 1)
    $headerCellB = $ws.Cells.Item($rowHeader, 2)
    $markerText = $headerCellA.Formula
    $regionText = $headerCellB.Formula

    # Row N: move the "~UC_Sets: T_xx: " label from B to A.
    $ws.Cells.Item($rowTitle, 1).Formula = $titleText
    $titleCellB.ClearContents() | Out-Null

    # Row N+1: swap A and B.
    $headerCellA.Formula = $regionText
    $headerCellB.Formula = $markerText
}

foreach ($ws in $wb.Worksheets) {
    # Find the "used" extent of the sheet so we know how far to scan.
    $lastRow = $ws.Cells.SpecialCells(11).Row   # xlCellTypeLastCell = 11

    for ($r = 1; $r -le $lastRow; $r++) {
        $bText = $ws.Cells.Item($r, 2).Formula
        if ($bText -like "~UC_Sets: T_*") {
            Swap-UcSetsBlock $ws $r
        }
    }
}
